$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Department column (C) split from the single generic "FACULTY OF
# HOSPITALITY" label into "Hospitality" for the single-qualification rows
# (2-4) and "Packages" for the combined-qualification rows (5-7).
$ws.Range("C2:C4").Value = "Hospitality"
$ws.Range("C5:C7").Value = "Packages"
